# Update the cryptocurrency price list with refreshed Price (D) and Volume(1h) (E) values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => new Price text (D column). $null means the Price column is unchanged for that row.
$priceUpdates = @(
    @{ Row = 2;  Value = "26.944.78" },
    @{ Row = 3;  Value = "1.555.04" },
    @{ Row = 5;  Value = "206.96" },
    @{ Row = 8;  Value = "22.09" },
    @{ Row = 10; Value = "0.0589" },
    @{ Row = 11; Value = "0.0858" },
    @{ Row = 12; Value = "1.775.99" },
    @{ Row = 13; Value = "1.553.66" },
    @{ Row = 15; Value = "0.520" },
    @{ Row = 16; Value = "26.929.77" },
    @{ Row = 17; Value = "61.78" },
    @{ Row = 18; Value = "218.05" },
    @{ Row = 19; Value = "0.0₃0697" },
    @{ Row = 20; Value = "7.32" },
    @{ Row = 23; Value = "9.22" },
    @{ Row = 24; Value = "1.95" },
    @{ Row = 25; Value = "154.10" },
    @{ Row = 26; Value = "6.63" },
    @{ Row = 27; Value = "14.95" },
    @{ Row = 31; Value = "1.09" },
    @{ Row = 32; Value = "3.22" },
    @{ Row = 33; Value = "1.427.46" },
    @{ Row = 36; Value = "0.975" },
    @{ Row = 37; Value = "2.28" },
    @{ Row = 38; Value = "0.0164" },
    @{ Row = 39; Value = "0.520" },
    @{ Row = 40; Value = "0.814" },
    @{ Row = 42; Value = "5.70" },
    @{ Row = 43; Value = "2.29" },
    @{ Row = 44; Value = "0.985" },
    @{ Row = 45; Value = "64.55" },
    @{ Row = 46; Value = "1.75" },
    @{ Row = 47; Value = "1.689.82" },
    @{ Row = 48; Value = "87.60" },
    @{ Row = 50; Value = "0.0₆0100" }
)

# Row => new Volume(1h) text (E column, with literal leading/trailing double spaces).
$volumeUpdates = @(
    @{ Row = 2;  Value = "  +0.17%  " },
    @{ Row = 3;  Value = "  +0.46%  " },
    @{ Row = 4;  Value = "  -0.11%  " },
    @{ Row = 5;  Value = "  +0.03%  " },
    @{ Row = 6;  Value = "  +0.04%  " },
    @{ Row = 7;  Value = "  -0.15%  " },
    @{ Row = 8;  Value = "  +2.98%  " },
    @{ Row = 9;  Value = "  +0.34%  " },
    @{ Row = 10; Value = "  +1.00%  " },
    @{ Row = 11; Value = "  +0.35%  " },
    @{ Row = 12; Value = "  +0.46%  " },
    @{ Row = 13; Value = "  +0.35%  " },
    @{ Row = 14; Value = "  +1.11%  " },
    @{ Row = 15; Value = "  +1.60%  " },
    @{ Row = 16; Value = "  +0.09%  " },
    @{ Row = 17; Value = "  +0.32%  " },
    @{ Row = 18; Value = "  +1.44%  " },
    @{ Row = 19; Value = "  +2.00%  " },
    @{ Row = 20; Value = "  +1.14%  " },
    @{ Row = 21; Value = "  -0.04%  " },
    @{ Row = 22; Value = "  +1.28%  " },
    @{ Row = 23; Value = "  +0.35%  " },
    @{ Row = 24; Value = "  +0.60%  " },
    @{ Row = 25; Value = "  +1.29%  " },
    @{ Row = 26; Value = "  -0.14%  " },
    @{ Row = 27; Value = "  +0.40%  " },
    @{ Row = 28; Value = "  +0.72%  " },
    @{ Row = 29; Value = "  -0.13%  " },
    @{ Row = 30; Value = "  +2.15%  " },
    @{ Row = 31; Value = "  -0.58%  " },
    @{ Row = 32; Value = "  -0.22%  " },
    @{ Row = 33; Value = "  +4.21%  " },
    @{ Row = 34; Value = "  +4.36%  " },
    @{ Row = 35; Value = "  +3.23%  " },
    @{ Row = 36; Value = "  +1.99%  " },
    @{ Row = 37; Value = "  +0.04%  " },
    @{ Row = 38; Value = "  -0.14%  " },
    @{ Row = 39; Value = "  -0.20%  " },
    @{ Row = 40; Value = "  +0.73%  " },
    @{ Row = 41; Value = "  -0.12%  " },
    @{ Row = 42; Value = "  +1.04%  " },
    @{ Row = 43; Value = "  +2.94%  " },
    @{ Row = 44; Value = "  -0.43%  " },
    @{ Row = 45; Value = "  +1.52%  " },
    @{ Row = 46; Value = "  +1.08%  " },
    @{ Row = 47; Value = "  +0.42%  " },
    @{ Row = 48; Value = "  +2.69%  " },
    @{ Row = 49; Value = "  +2.54%  " },
    @{ Row = 50; Value = "  +3.12%  " },
    @{ Row = 51; Value = "  +1.38%  " }
)

foreach ($u in $priceUpdates) {
    $cell = $ws.Cells.Item($u.Row, 4)
    # Prefix with an apostrophe so Excel always stores the value as text, then
    # reset the style back to Normal so no stray "quote prefix" text format
    # (and its associated style index) gets attached to the cell.
    $cell.Value = "'" + $u.Value
    $cell.Style = "Normal"
}

foreach ($u in $volumeUpdates) {
    $ws.Cells.Item($u.Row, 5).Value = $u.Value
}
